$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.55 = 17954.55 pesos`n✅ 17954.55 pesos = 4.52 = 958.3 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 220
$ws2.Range("O10").Value = 3950
$ws2.Range("N12").Value = 3972
$ws2.Range("O12").Value = 212
